$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- First-page header (header1.xml): BTec logo "image1.jpg" -> "image2.jpg" ---
$hdrFirst = $sec.Headers.Item(2)
$hdrShape = $hdrFirst.Range.InlineShapes.Item(1)
$hdrShape.Range.InlineShapes.Item(1).Name = "image2.jpg"

# --- First-page footer (footer1.xml): Pearson logo "image2.png" -> "image1.png" ---
$ftrFirst = $sec.Footers.Item(2)
$ftrFirstShape = $ftrFirst.Range.InlineShapes.Item(1)
$ftrFirstShape.Range.InlineShapes.Item(1).Name = "image1.png"

# --- Default/primary footer (footer2.xml): Pearson logo "image2.png" -> "image1.png" ---
$ftrDefault = $sec.Footers.Item(1)
$ftrDefaultShape = $ftrDefault.Range.InlineShapes.Item(1)
$ftrDefaultShape.Range.InlineShapes.Item(1).Name = "image1.png"

Write-Host "Header(2) shape name:" $hdrFirst.Range.InlineShapes.Item(1).Name
Write-Host "Footer(2) shape name:" $ftrFirst.Range.InlineShapes.Item(1).Name
Write-Host "Footer(1) shape name:" $ftrDefault.Range.InlineShapes.Item(1).Name
